$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.877.14"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.546.01"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'567.45"
$ws.Range("D6").Value = "'147.08"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "2.541.02"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "'5.60"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'27.24"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").Value = "2.997.64"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("D16").Value = "62.835.36"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "2.541.60"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "'11.44"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "'336.73"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "'4.29"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'65.29"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'1.63"
$ws.Range("E25").Value = "  +8.81%  "
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").Value = "'1.51"
$ws.Range("E27").Value = "  +10.98%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "'8.38"
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "'7.28"
$ws.Range("E30").Value = "  +6.55%  "
$ws.Range("D31").Value = "0.0₃0812"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'178.13"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'407.00"
$ws.Range("E35").Value = "  +10.14%  "
$ws.Range("D36").Value = "'0.399"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "'18.94"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D39").Value = "'4.35"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'39.10"
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("D43").Value = "'153.01"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").Value = "'3.74"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "'20.78"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "'0.606"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'0.0959"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'0.0517"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").Value = "'18.20"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("E51").Value = "  +0.50%  "
